$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.091.28"
$ws.Range("E2").Value = "  +3.87%  "
$ws.Range("D3").Value = "4.034.60"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'519.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").Value = "'146.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("D7").Value = "'0.734"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.03%  "
$ws.Range("D8").Value = "4.028.85"
$ws.Range("E8").Value = "  +3.47%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("D10").Value = "'0.777"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.94%  "
$ws.Range("D11").Value = "'0.175"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'0.0000328"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").Value = "'47.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.08%  "
$ws.Range("D14").Value = "'11.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.98%  "
$ws.Range("D15").Value = "4.678.25"
$ws.Range("E15").Value = "  +3.48%  "
$ws.Range("D16").Value = "4.035.93"
$ws.Range("E16").Value = "  +3.70%  "
$ws.Range("E17").Value = "  +6.85%  "
$ws.Range("D18").Value = "'14.17"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  -0.76%  "
$ws.Range("D21").Value = "71.980.39"
$ws.Range("E21").Value = "  +3.75%  "
$ws.Range("D22").Value = "'444.09"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.41%  "
$ws.Range("D23").Value = "'104.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +19.07%  "
$ws.Range("D24").Value = "'3.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.78%  "
$ws.Range("D25").Value = "'14.84"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.85%  "
$ws.Range("D26").Value = "'4.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.64%  "
$ws.Range("D27").Value = "'11.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "
$ws.Range("D28").Value = "'11.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.04%  "
$ws.Range("D29").Value = "'37.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.85%  "
$ws.Range("D30").Value = "'5.80"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.23%  "
$ws.Range("D31").Value = "'3.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.64%  "
$ws.Range("D32").Value = "'13.66"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("E33").Value = "  +3.02%  "
$ws.Range("D34").Value = "'682.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.92%  "
$ws.Range("D35").Value = "'6.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +14.77%  "
$ws.Range("D36").Value = "'66.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.58%  "
$ws.Range("D37").Value = "'42.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.71%  "
$ws.Range("E38").Value = "  -1.02%  "
$ws.Range("E39").Value = "  -1.49%  "
$ws.Range("D40").Value = "'3.54"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.22%  "
$ws.Range("D41").Value = "'0.151"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  +3.30%  "
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").Value = "'0.160"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +13.92%  "
$ws.Range("D47").Value = "'3.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.84%  "
$ws.Range("E48").Value = "  -2.99%  "
$ws.Range("D49").Value = "'3.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("D50").Value = "'9.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.70%  "
$ws.Range("D51").Value = "'3.35"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.62%  "
